$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format column A (dates) as text before writing so Excel does not auto-convert
# dd/mm/yyyy strings into date serial numbers.
$ws.Range("A3:A45").NumberFormat = "@"

$ws.Range("A1").Value = "Tanggal"
$ws.Range("B1").Value = "Keterangan"
$ws.Range("C1").Value = "Debit"
$ws.Range("D1").Value = "Kredit"
$ws.Range("E1").Value = "Saldo"

$ws.Range("A2").Value = $null
$ws.Range("B2").Value = "SALDO AWAL"
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 318599650.03

$ws.Range("A3").Value = "01/01/2025"
$ws.Range("B3").Value = "MID:001770398"
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 13483175
$ws.Range("E3").Value = 332082825.03

$ws.Range("A4").Value = "02/01/2025"
$ws.Range("B4").Value = "MID:001770398"
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 6525848
$ws.Range("E4").Value = 0

$ws.Range("A5").Value = "02/01/2025"
$ws.Range("B5").Value = "O2O1/FTSCY/WS95051"
$ws.Range("C5").Value = 27050000
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 311558673.03

$ws.Range("A6").Value = "03/01/2025"
$ws.Range("B6").Value = "MID:001770398"
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 33265241
$ws.Range("E6").Value = 344823914.03

$ws.Range("A7").Value = "04/01/2025"
$ws.Range("B7").Value = "MID:001770398"
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 9423862
$ws.Range("E7").Value = 354247776.03

$ws.Range("A8").Value = "05/01/2025"
$ws.Range("B8").Value = "MID:001770398"
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 8572629
$ws.Range("E8").Value = 0

$ws.Range("A9").Value = "05/01/2025"
$ws.Range("B9").Value = "MID 885001770398"
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 1023153.12
$ws.Range("E9").Value = 363843558.15

$ws.Range("A10").Value = "06/01/2025"
$ws.Range("B10").Value = "MID:001770398"
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 8164244
$ws.Range("E10").Value = 372007802.15

$ws.Range("A11").Value = "07/01/2025"
$ws.Range("B11").Value = "O6O1/FTSCY/WS95051"
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 6121898
$ws.Range("E11").Value = 378129700.15

$ws.Range("A12").Value = "07/01/2025"
$ws.Range("B12").Value = "MID:001770398"
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 7306984
$ws.Range("E12").Value = 385436684.15

$ws.Range("A13").Value = "08/01/2025"
$ws.Range("B13").Value = "MID:001770398"
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 5570888
$ws.Range("E13").Value = 391007572.15

$ws.Range("A14").Value = "09/01/2025"
$ws.Range("B14").Value = "MID:001770398"
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 4360016
$ws.Range("E14").Value = 0

$ws.Range("A15").Value = "09/01/2025"
$ws.Range("B15").Value = "O9O1/FTSCY/WS95051"
$ws.Range("C15").Value = 244730000
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 150637588.15

$ws.Range("A16").Value = "10/01/2025"
$ws.Range("B16").Value = "MID:001770398"
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 9181621
$ws.Range("E16").Value = 159819209.15

$ws.Range("A17").Value = "11/01/2025"
$ws.Range("B17").Value = "MID:001770398"
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 20261504
$ws.Range("E17").Value = 180080713.15

$ws.Range("A18").Value = "12/01/2025"
$ws.Range("B18").Value = "MID:001770398"
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 10643898
$ws.Range("E18").Value = 190724611.15

$ws.Range("A19").Value = "13/01/2025"
$ws.Range("B19").Value = "MID:001770398"
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 11136960
$ws.Range("E19").Value = 201861571.15

$ws.Range("A20").Value = "14/01/2025"
$ws.Range("B20").Value = "MID:001770398"
$ws.Range("C20").Value = 0
$ws.Range("D20").Value = 2628790
$ws.Range("E20").Value = 204490361.15

$ws.Range("A21").Value = "14/01/2025"
$ws.Range("B21").Value = "14O1/FTSCY/WS95051"
$ws.Range("C21").Value = 35550000
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 168940361.15

$ws.Range("A22").Value = "15/01/2025"
$ws.Range("B22").Value = "MID:001770398"
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 5139229
$ws.Range("E22").Value = 0

$ws.Range("A23").Value = "15/01/2025"
$ws.Range("B23").Value = "NTRF@1185719261"
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 6875000
$ws.Range("E23").Value = 180954590.15

$ws.Range("A24").Value = "16/01/2025"
$ws.Range("B24").Value = "MID:001770398"
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 7728376
$ws.Range("E24").Value = 188682966.15

$ws.Range("A25").Value = "17/01/2025"
$ws.Range("B25").Value = "MID:001770398"
$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 6895087
$ws.Range("E25").Value = 0

$ws.Range("A26").Value = "17/01/2025"
$ws.Range("B26").Value = $null
$ws.Range("C26").Value = 20000
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = 195558053.15

$ws.Range("A27").Value = "18/01/2025"
$ws.Range("B27").Value = "MID:001770398"
$ws.Range("C27").Value = 0
$ws.Range("D27").Value = 8014683
$ws.Range("E27").Value = 203572736.15

$ws.Range("A28").Value = "19/01/2025"
$ws.Range("B28").Value = "MID:001770398"
$ws.Range("C28").Value = 0
$ws.Range("D28").Value = 6255144
$ws.Range("E28").Value = 209827880.15

$ws.Range("A29").Value = "20/01/2025"
$ws.Range("B29").Value = "MID:001770398"
$ws.Range("C29").Value = 0
$ws.Range("D29").Value = 10723347
$ws.Range("E29").Value = 220551227.15

$ws.Range("A30").Value = "21/01/2025"
$ws.Range("B30").Value = "MID:001770398"
$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 2563183
$ws.Range("E30").Value = 223114410.15

$ws.Range("A31").Value = "21/01/2025"
$ws.Range("B31").Value = "DR 028"
$ws.Range("C31").Value = 0
$ws.Range("D31").Value = 925000
$ws.Range("E31").Value = 224039410.15

$ws.Range("A32").Value = "22/01/2025"
$ws.Range("B32").Value = "MID:001770398"
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 17828643
$ws.Range("E32").Value = 0

$ws.Range("A33").Value = "22/01/2025"
$ws.Range("B33").Value = "22O1/FTSCY/WS95051"
$ws.Range("C33").Value = 11125000
$ws.Range("D33").Value = 0
$ws.Range("E33").Value = 230743053.15

$ws.Range("A34").Value = "23/01/2025"
$ws.Range("B34").Value = "MID:001770398"
$ws.Range("C34").Value = 0
$ws.Range("D34").Value = 4094940
$ws.Range("E34").Value = 234837993.15

$ws.Range("A35").Value = "24/01/2025"
$ws.Range("B35").Value = "MID:001770398"
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = 4968134
$ws.Range("E35").Value = 239806127.15

$ws.Range("A36").Value = "25/01/2025"
$ws.Range("B36").Value = "MID:001770398"
$ws.Range("C36").Value = 0
$ws.Range("D36").Value = 7988959
$ws.Range("E36").Value = 247795086.15

$ws.Range("A37").Value = "26/01/2025"
$ws.Range("B37").Value = "MID:001770398"
$ws.Range("C37").Value = 0
$ws.Range("D37").Value = 8718225
$ws.Range("E37").Value = 256513311.15

$ws.Range("A38").Value = "27/01/2025"
$ws.Range("B38").Value = "MID:001770398"
$ws.Range("C38").Value = 0
$ws.Range("D38").Value = 14955267
$ws.Range("E38").Value = 271468578.15

$ws.Range("A39").Value = "28/01/2025"
$ws.Range("B39").Value = "MID:001770398"
$ws.Range("C39").Value = 0
$ws.Range("D39").Value = 8589883
$ws.Range("E39").Value = 280058461.15

$ws.Range("A40").Value = "29/01/2025"
$ws.Range("B40").Value = "MID:001770398"
$ws.Range("C40").Value = 0
$ws.Range("D40").Value = 2772665
$ws.Range("E40").Value = 282831126.15

$ws.Range("A41").Value = "29/01/2025"
$ws.Range("B41").Value = "29O1/FTSCY/WS95051"
$ws.Range("C41").Value = 66415000
$ws.Range("D41").Value = 0
$ws.Range("E41").Value = 216416126.15

$ws.Range("A42").Value = "30/01/2025"
$ws.Range("B42").Value = "MID:001770398"
$ws.Range("C42").Value = 0
$ws.Range("D42").Value = 1767600
$ws.Range("E42").Value = 218183726.15

$ws.Range("A43").Value = "31/01/2025"
$ws.Range("B43").Value = "MID:001770398"
$ws.Range("C43").Value = 0
$ws.Range("D43").Value = 5317589
$ws.Range("E43").Value = 0

$ws.Range("A44").Value = "31/01/2025"
$ws.Range("B44").Value = $null
$ws.Range("C44").Value = 0
$ws.Range("D44").Value = 2125.07
$ws.Range("E44").Value = 0

$ws.Range("A45").Value = "31/01/2025"
$ws.Range("B45").Value = $null
$ws.Range("C45").Value = 425.01
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 223503015.21

$ws.Range("A46").Value = $null
$ws.Range("B46").Value = "SALDO AKHIR"
$ws.Range("C46").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 223503015.21

Write-Host "Updated sheet data through row 46"
